$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("B2").Value = 4
$ws1.Range("D2").Value = 8.699999999999999
$ws1.Range("E2").Value = 1.3

$ws1.Range("B3").Value = 7
$ws1.Range("C3").Value = 8
$ws1.Range("D3").Value = 15.2
$ws1.Range("E3").Value = 10.1

$ws1.Range("C4").Value = 63
$ws1.Range("D4").Value = 56.5
$ws1.Range("E4").Value = 79.7

$ws1.Range("C5").Value = 5
$ws1.Range("D5").Value = 10.9
$ws1.Range("E5").Value = 6.3

$ws1.Range("B6").Value = 4
$ws1.Range("C6").Value = 2
$ws1.Range("D6").Value = 8.699999999999999
$ws1.Range("E6").Value = 2.5

$ws1.Range("B7").Value = 155
$ws1.Range("C7").Value = 203

# --- Sheet: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("C3").Value = 46
$ws4.Range("C4").Value = 79

# --- Sheet: "Interannual update - High Pri" ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")

$ws5.Range("B2").Value = 62
$ws5.Range("C2").Value = 60.2

$ws5.Range("B4").Value = 40
$ws5.Range("C4").Value = 38.8
